$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4716.2
$ws.Range("I43").Value = 9500.5
$ws.Range("J43").Value = 1526.6666
$ws.Range("K43").Value = 9500.5
$ws.Range("L43").Value = 1526.6666
$ws.Range("M43").Value = -9431.5
$ws.Range("N43").Value = -1664.6666

$ws.Range("H114").Value = 24788.8
$ws.Range("J114").Value = 24788.8
$ws.Range("L114").Value = 24788.8
$ws.Range("N114").Value = -33466.8

$ws.Range("H123").Value = 41320
$ws.Range("J123").Value = 41320
$ws.Range("L123").Value = 41320
$ws.Range("N123").Value = -51120

$ws.Range("H127").Value = 2252.3171
$ws.Range("I127").Value = 876.7778
$ws.Range("K127").Value = 2630.3334
$ws.Range("M127").Value = 2329.6666

$ws.Range("H128").Value = 35920
$ws.Range("J128").Value = 35920
$ws.Range("L128").Value = 35920
$ws.Range("N128").Value = -45880

$ws.Range("H138").Value = 4082.75
$ws.Range("I138").Value = 891.4706
$ws.Range("J138").Value = 4736.3857
$ws.Range("K138").Value = 2674.4118
$ws.Range("L138").Value = 14209.1571
$ws.Range("M138").Value = 2465.5882
$ws.Range("N138").Value = -24489.1571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2420.75
$ws.Range("I2").Value = 2753.8333
$ws.Range("J2").Value = 2087.6667
$ws.Range("K2").Value = 2753.8333
$ws.Range("L2").Value = 2087.6667
$ws.Range("M2").Value = -2640.8333
$ws.Range("N2").Value = -2313.6667

$ws.Range("H45").Value = 1852.8572
$ws.Range("I45").Value = 1698.4546
$ws.Range("J45").Value = 2419
$ws.Range("K45").Value = 1698.4546
$ws.Range("L45").Value = 2419
$ws.Range("M45").Value = -1321.4546
$ws.Range("N45").Value = -3173

$ws.Range("H110").Value = 599.7143
$ws.Range("I110").Value = 582
$ws.Range("J110").Value = 706
$ws.Range("K110").Value = 582
$ws.Range("L110").Value = 706
$ws.Range("M110").Value = 1463
$ws.Range("N110").Value = -4796

$ws.Range("H116").Value = 2420.75
$ws.Range("I116").Value = 2753.8333
$ws.Range("J116").Value = 2087.6667
$ws.Range("K116").Value = 2753.8333
$ws.Range("L116").Value = 2087.6667
$ws.Range("M116").Value = -459.8332999999998
$ws.Range("N116").Value = -6675.6667

$ws.Range("H132").Value = 1657.7894
$ws.Range("I132").Value = 1031.3125
$ws.Range("K132").Value = 3093.9375
$ws.Range("M132").Value = -563.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2420.75
$ws.Range("I3").Value = 2753.8333
$ws.Range("J3").Value = 2087.6667
$ws.Range("K3").Value = 2753.8333
$ws.Range("L3").Value = 2087.6667
$ws.Range("M3").Value = -2639.8333
$ws.Range("N3").Value = -2315.6667

$ws.Range("H20").Value = 40025504
$ws.Range("I20").Value = 44223.43
$ws.Range("J20").Value = 90910776
$ws.Range("K20").Value = 44223.43
$ws.Range("L20").Value = 90910776
$ws.Range("M20").Value = -43976.43
$ws.Range("N20").Value = -90911270

$ws.Range("H42").Value = 120000
$ws.Range("J42").Value = 120000
$ws.Range("L42").Value = 120000
$ws.Range("N42").Value = -120656

$ws.Range("H94").Value = 976.15
$ws.Range("I94").Value = 908.2
$ws.Range("K94").Value = 908.2
$ws.Range("M94").Value = -457.2

$ws.Range("H105").Value = 8004.3213
$ws.Range("I105").Value = 7582.778
$ws.Range("J105").Value = 8763.1
$ws.Range("K105").Value = 7582.778
$ws.Range("L105").Value = 8763.1
$ws.Range("M105").Value = -5835.778
$ws.Range("N105").Value = -12257.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2005.0613
$ws.Range("I99").Value = 1582.2258
$ws.Range("J99").Value = 2733.2778
$ws.Range("K99").Value = 1582.2258
$ws.Range("L99").Value = 2733.2778
$ws.Range("M99").Value = -84.22579999999994
$ws.Range("N99").Value = -5729.2778

$ws.Range("H126").Value = 2005.0613
$ws.Range("I126").Value = 1582.2258
$ws.Range("J126").Value = 2733.2778
$ws.Range("K126").Value = 4746.6774
$ws.Range("L126").Value = 8199.8334
$ws.Range("M126").Value = -2276.6774
$ws.Range("N126").Value = -13139.8334

$ws.Range("H141").Value = 48712.5
$ws.Range("J141").Value = 48712.5
$ws.Range("L141").Value = 48712.5
$ws.Range("N141").Value = -59072.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1173.4445
$ws.Range("I5").Value = 831.7727
$ws.Range("J5").Value = 1710.3572
$ws.Range("K5").Value = 2495.3181
$ws.Range("L5").Value = 5131.071599999999
$ws.Range("M5").Value = -2383.3181
$ws.Range("N5").Value = -5355.071599999999

$ws.Range("H107").Value = 556412.0600000001
$ws.Range("J107").Value = 1556134.2
$ws.Range("L107").Value = 4668402.6
$ws.Range("N107").Value = -4672242.6

$ws.Range("H131").Value = 8948186
$ws.Range("I131").Value = 55668140
$ws.Range("J131").Value = 1811.9362
$ws.Range("K131").Value = 167004420
$ws.Range("L131").Value = 5435.8086
$ws.Range("M131").Value = -166999380
$ws.Range("N131").Value = -15515.8086

$ws.Range("H135").Value = 1173.4445
$ws.Range("I135").Value = 831.7727
$ws.Range("J135").Value = 1710.3572
$ws.Range("K135").Value = 7485.954299999999
$ws.Range("L135").Value = 15393.2148
$ws.Range("M135").Value = -4950.954299999999
$ws.Range("N135").Value = -20463.2148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4228.5713
$ws.Range("I113").Value = 6110.1113
$ws.Range("J113").Value = 841.8
$ws.Range("K113").Value = 6110.1113
$ws.Range("L113").Value = 841.8
$ws.Range("M113").Value = -3940.1113
$ws.Range("N113").Value = -5181.8

$ws.Range("H126").Value = 2665
$ws.Range("I126").Value = 3622.4
$ws.Range("J126").Value = 1867.1666
$ws.Range("K126").Value = 10867.2
$ws.Range("L126").Value = 5601.4998
$ws.Range("M126").Value = -8397.200000000001
$ws.Range("N126").Value = -10541.4998

$ws.Range("H132").Value = 2680.0303
$ws.Range("I132").Value = 2507.4348
$ws.Range("J132").Value = 3077
$ws.Range("K132").Value = 7522.3044
$ws.Range("L132").Value = 9231
$ws.Range("M132").Value = -4992.3044
$ws.Range("N132").Value = -14291

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3113
$ws.Range("I93").Value = 2000
$ws.Range("J93").Value = 3272
$ws.Range("K93").Value = 2000
$ws.Range("L93").Value = 3272
$ws.Range("M93").Value = -752
$ws.Range("N93").Value = -5768

$ws.Range("H132").Value = 1845.9642
$ws.Range("I132").Value = 1322.2941
$ws.Range("J132").Value = 2655.2727
$ws.Range("K132").Value = 3966.8823
$ws.Range("L132").Value = 7965.8181
$ws.Range("M132").Value = -1436.8823
$ws.Range("N132").Value = -13025.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2095.7058
$ws.Range("I81").Value = 1706
$ws.Range("J81").Value = 2442.111
$ws.Range("K81").Value = 3412
$ws.Range("L81").Value = 4884.222
$ws.Range("M81").Value = -2351
$ws.Range("N81").Value = -7006.222

$ws.Range("H84").Value = 2095.7058
$ws.Range("I84").Value = 1706
$ws.Range("J84").Value = 2442.111
$ws.Range("K84").Value = 17060
$ws.Range("L84").Value = 24421.11
$ws.Range("M84").Value = -11756
$ws.Range("N84").Value = -35029.11

$ws.Range("H96").Value = 3600
$ws.Range("I96").Value = 3600
$ws.Range("K96").Value = 3600
$ws.Range("M96").Value = -2227

$ws.Range("H107").Value = 533.75
$ws.Range("I107").Value = 527.3
$ws.Range("J107").Value = 544.5
$ws.Range("K107").Value = 1581.9
$ws.Range("L107").Value = 1633.5
$ws.Range("M107").Value = 338.1000000000001
$ws.Range("N107").Value = -5473.5

$ws.Range("H122").Value = 2413.1428
$ws.Range("I122").Value = 2353.7778
$ws.Range("J122").Value = 2520
$ws.Range("K122").Value = 7061.3334
$ws.Range("L122").Value = 7560
$ws.Range("M122").Value = -4611.3334
$ws.Range("N122").Value = -12460

$ws.Range("H126").Value = 817.6667
$ws.Range("I126").Value = 818.3333
$ws.Range("J126").Value = 816.6667
$ws.Range("K126").Value = 2454.9999
$ws.Range("L126").Value = 2450.0001
$ws.Range("M126").Value = 15.0001000000002
$ws.Range("N126").Value = -7390.0001

$ws.Range("H132").Value = 1500.2424
$ws.Range("I132").Value = 832.2727
$ws.Range("J132").Value = 2836.182
$ws.Range("K132").Value = 2496.8181
$ws.Range("L132").Value = 8508.545999999998
$ws.Range("M132").Value = 33.18190000000004
$ws.Range("N132").Value = -13568.546
